$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2023
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E2").Value = 901745.342027444
$ws.Range("F2").Value = 755558432.84
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2023
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "Servicios de limpieza industrial"
$ws.Range("E3").Value = 619549.5905107583
$ws.Range("F3").Value = 510549401.65
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 2

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2023
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "Equipamiento y suministros médicos"
$ws.Range("E4").Value = 279168.0421732903
$ws.Range("F4").Value = 232866300.68
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 3

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2023
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E5").Value = 252417.4300486843
$ws.Range("F5").Value = 212712186
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 4

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2023
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "Artículos eléctricos y de iluminación"
$ws.Range("E6").Value = 185353.64600234
$ws.Range("F6").Value = 162446749
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 5

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 2023
$ws.Range("C7").Value = "Antofagasta"
$ws.Range("D7").Value = "Salud, servicios sanitarios y alimentación"
$ws.Range("E7").Value = 200782582.7296285
$ws.Range("F7").Value = 172210606470.1
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2023
$ws.Range("C8").Value = "Antofagasta"
$ws.Range("D8").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E8").Value = 43525675.8456754
$ws.Range("F8").Value = 37161120531.35535
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 2

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 2023
$ws.Range("C9").Value = "Antofagasta"
$ws.Range("D9").Value = "Servicios de limpieza industrial"
$ws.Range("E9").Value = 28528175.3848198
$ws.Range("F9").Value = 24130276882.70694
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 3

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 2023
$ws.Range("C10").Value = "Antofagasta"
$ws.Range("D10").Value = "Obras MINVU"
$ws.Range("E10").Value = 25665374.41055844
$ws.Range("F10").Value = 20807430211
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 4

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 2023
$ws.Range("C11").Value = "Antofagasta"
$ws.Range("D11").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E11").Value = 22941684.16735266
$ws.Range("F11").Value = 19330580572.54248
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 5

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 2023
$ws.Range("C12").Value = "Araucanía"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = 116133915.7099773
$ws.Range("F12").Value = 96705116160
$ws.Range("G12").Value = 2681005.769840984
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 1

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 2023
$ws.Range("C13").Value = "Araucanía"
$ws.Range("D13").Value = "Equipamiento y suministros médicos"
$ws.Range("E13").Value = 64125357.10944952
$ws.Range("F13").Value = 53517744694.78053
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 2

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 2023
$ws.Range("C14").Value = "Araucanía"
$ws.Range("D14").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E14").Value = 63787457.90708901
$ws.Range("F14").Value = 53485522347.2868
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 3

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 2023
$ws.Range("C15").Value = "Araucanía"
$ws.Range("D15").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E15").Value = 54422704.98659359
$ws.Range("F15").Value = 45855734441.98038
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 4

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 2023
$ws.Range("C16").Value = "Araucanía"
$ws.Range("D16").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E16").Value = 45258967.52201545
$ws.Range("F16").Value = 37890607898.8496
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 5

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 2023
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = 1352270603090615
$ws.Range("F17").Value = 1087004945419016960
$ws.Range("G17").Value = 30476785573881.11
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 1

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = "Servicios basados en ingeniería, ciencias sociales y tecnología de la información"
$ws.Range("E18").Value = 35654089.69784889
$ws.Range("F18").Value = 30113570613.80602
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 2

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 2023
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E19").Value = 21445120.03490136
$ws.Range("F19").Value = 17842193140.04343
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 3

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 2023
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = "Servicios de limpieza industrial"
$ws.Range("E20").Value = 20099303.14632814
$ws.Range("F20").Value = 16413677617.78473
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 4

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 2023
$ws.Range("C21").Value = "Arica y Parinacota"
$ws.Range("D21").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E21").Value = 19330050.16764392
$ws.Range("F21").Value = 16093992591.11241
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 5

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 2023
$ws.Range("C22").Value = "Atacama"
$ws.Range("D22").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E22").Value = 54531112.75793689
$ws.Range("F22").Value = 45862515884.4021
$ws.Range("G22").Value = ""
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 1

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 2023
$ws.Range("C23").Value = "Atacama"
$ws.Range("D23").Value = "Equipamiento y suministros médicos"
$ws.Range("E23").Value = 13068163.0374981
$ws.Range("F23").Value = 10925217940.37981
$ws.Range("G23").Value = ""
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 2

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 2023
$ws.Range("C24").Value = "Atacama"
$ws.Range("D24").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E24").Value = 12405223.91676243
$ws.Range("F24").Value = 10611998519.18348
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 3

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 2023
$ws.Range("C25").Value = "Atacama"
$ws.Range("D25").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E25").Value = 11604608.89120818
$ws.Range("F25").Value = 9673860675.228569
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 4

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 2023
$ws.Range("C26").Value = "Atacama"
$ws.Range("D26").Value = ""
$ws.Range("E26").Value = 8108010.321714605
$ws.Range("F26").Value = 6923824055
$ws.Range("G26").Value = 191962.0064127275
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 5

$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 2023
$ws.Range("C27").Value = "Aysén"
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = 75181925535215.67
$ws.Range("F27").Value = 62664496549964416
$ws.Range("G27").Value = 1735822490109.256
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 1

$ws.Range("A28").Value = 26
$ws.Range("B28").Value = 2023
$ws.Range("C28").Value = "Aysén"
$ws.Range("D28").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E28").Value = 31177870.90745586
$ws.Range("F28").Value = 26521149604.0243
$ws.Range("G28").Value = ""
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 2

$ws.Range("A29").Value = 27
$ws.Range("B29").Value = 2023
$ws.Range("C29").Value = "Aysén"
$ws.Range("D29").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E29").Value = 14003763.78279521
$ws.Range("F29").Value = 11775102934.55973
$ws.Range("G29").Value = ""
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 3

$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 2023
$ws.Range("C30").Value = "Aysén"
$ws.Range("D30").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E30").Value = 11009973.86117144
$ws.Range("F30").Value = 9294968316.416376
$ws.Range("G30").Value = ""
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 4

$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 2023
$ws.Range("C31").Value = "Aysén"
$ws.Range("D31").Value = "Equipamiento y suministros médicos"
$ws.Range("E31").Value = 9536607.435378332
$ws.Range("F31").Value = 8064745659.00032
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 5

$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 2023
$ws.Range("C32").Value = "Bío-Bío"
$ws.Range("D32").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E32").Value = 116062214.9149046
$ws.Range("F32").Value = 98668785006.57314
$ws.Range("G32").Value = ""
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 1

$ws.Range("A33").Value = 31
$ws.Range("B33").Value = 2023
$ws.Range("C33").Value = "Bío-Bío"
$ws.Range("D33").Value = "Equipamiento y suministros médicos"
$ws.Range("E33").Value = 110198246.1385839
$ws.Range("F33").Value = 92384853982.52975
$ws.Range("G33").Value = ""
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 2

$ws.Range("A34").Value = 32
$ws.Range("B34").Value = 2023
$ws.Range("C34").Value = "Bío-Bío"
$ws.Range("D34").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E34").Value = 55556721.0422359
$ws.Range("F34").Value = 46340282395.06536
$ws.Range("G34").Value = ""
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 3

$ws.Range("A35").Value = 33
$ws.Range("B35").Value = 2023
$ws.Range("C35").Value = "Bío-Bío"
$ws.Range("D35").Value = ""
$ws.Range("E35").Value = 54240121.2510023
$ws.Range("F35").Value = 44404038953
$ws.Range("G35").Value = 1238042.772370678
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 4

$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 2023
$ws.Range("C36").Value = "Bío-Bío"
$ws.Range("D36").Value = "Servicios de limpieza industrial"
$ws.Range("E36").Value = 53342830.59164016
$ws.Range("F36").Value = 44132044935.6428
$ws.Range("G36").Value = ""
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 5

$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 2023
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = ""
$ws.Range("E37").Value = 297771937195863.6
$ws.Range("F37").Value = 260441208331486496
$ws.Range("G37").Value = 7102417245230.266
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 1

$ws.Range("A38").Value = 36
$ws.Range("B38").Value = 2023
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E38").Value = 41633397.43545135
$ws.Range("F38").Value = 35225171237.5635
$ws.Range("G38").Value = ""
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 2

$ws.Range("A39").Value = 37
$ws.Range("B39").Value = 2023
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = "Equipamiento y suministros médicos"
$ws.Range("E39").Value = 39688754.07025614
$ws.Range("F39").Value = 33178669559.64686
$ws.Range("G39").Value = ""
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 3

$ws.Range("A40").Value = 38
$ws.Range("B40").Value = 2023
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E40").Value = 22735563.36148095
$ws.Range("F40").Value = 19051480355.64533
$ws.Range("G40").Value = ""
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 4

$ws.Range("A41").Value = 39
$ws.Range("B41").Value = 2023
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = "Servicios de limpieza industrial"
$ws.Range("E41").Value = 21923195.17128803
$ws.Range("F41").Value = 18733871908.928
$ws.Range("G41").Value = ""
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 5

$ws.Range("A42").Value = 40
$ws.Range("B42").Value = 2023
$ws.Range("C42").Value = "Lib. Gral. Bdo. O'Higgins"
$ws.Range("D42").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E42").Value = 63196714.87987141
$ws.Range("F42").Value = 53210729461.0369
$ws.Range("G42").Value = ""
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 1

$ws.Range("A43").Value = 41
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = "Lib. Gral. Bdo. O'Higgins"
$ws.Range("D43").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E43").Value = 38511999.087263
$ws.Range("F43").Value = 32445528788.7043
$ws.Range("G43").Value = ""
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 2

$ws.Range("A44").Value = 42
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = "Lib. Gral. Bdo. O'Higgins"
$ws.Range("D44").Value = "Equipamiento y suministros médicos"
$ws.Range("E44").Value = 32081165.41547768
$ws.Range("F44").Value = 26841455077.30179
$ws.Range("G44").Value = ""
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 3

$ws.Range("A45").Value = 43
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = "Lib. Gral. Bdo. O'Higgins"
$ws.Range("D45").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E45").Value = 25397943.36118131
$ws.Range("F45").Value = 21207261273.46177
$ws.Range("G45").Value = ""
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 4

$ws.Range("A46").Value = 44
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = "Lib. Gral. Bdo. O'Higgins"
$ws.Range("D46").Value = "Servicios de limpieza industrial"
$ws.Range("E46").Value = 24393372.03552687
$ws.Range("F46").Value = 20252256439.76813
$ws.Range("G46").Value = ""
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 5

$ws.Range("A47").Value = 45
$ws.Range("B47").Value = 2023
$ws.Range("C47").Value = "Los Lagos"
$ws.Range("D47").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E47").Value = 123386252.6613988
$ws.Range("F47").Value = 103754722050.9581
$ws.Range("G47").Value = ""
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 1

$ws.Range("A48").Value = 46
$ws.Range("B48").Value = 2023
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = "Equipamiento y suministros médicos"
$ws.Range("E48").Value = 59098560.77217842
$ws.Range("F48").Value = 49579353199.74539
$ws.Range("G48").Value = ""
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 2

$ws.Range("A49").Value = 47
$ws.Range("B49").Value = 2023
$ws.Range("C49").Value = "Los Lagos"
$ws.Range("D49").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E49").Value = 30572914.56840733
$ws.Range("F49").Value = 25571536140.43476
$ws.Range("G49").Value = ""
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 3

$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 2023
$ws.Range("C50").Value = "Los Lagos"
$ws.Range("D50").Value = "Servicios de limpieza industrial"
$ws.Range("E50").Value = 25730540.22009587
$ws.Range("F50").Value = 21208058014.67456
$ws.Range("G50").Value = ""
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 4

$ws.Range("A51").Value = 49
$ws.Range("B51").Value = 2023
$ws.Range("C51").Value = "Los Lagos"
$ws.Range("D51").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E51").Value = 24915025.19560209
$ws.Range("F51").Value = 20783598106.66674
$ws.Range("G51").Value = ""
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 5

$ws.Range("A52").Value = 50
$ws.Range("B52").Value = 2023
$ws.Range("C52").Value = "Los Ríos"
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = 109617373.6895769
$ws.Range("F52").Value = 89537059870
$ws.Range("G52").Value = 2487471.73642611
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 1

$ws.Range("A53").Value = 51
$ws.Range("B53").Value = 2023
$ws.Range("C53").Value = "Los Ríos"
$ws.Range("D53").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E53").Value = 41389788.17639003
$ws.Range("F53").Value = 34863311480.4782
$ws.Range("G53").Value = ""
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 2

$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 2023
$ws.Range("C54").Value = "Los Ríos"
$ws.Range("D54").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E54").Value = 29345224.70572385
$ws.Range("F54").Value = 24827654181.91487
$ws.Range("G54").Value = ""
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 3

$ws.Range("A55").Value = 53
$ws.Range("B55").Value = 2023
$ws.Range("C55").Value = "Los Ríos"
$ws.Range("D55").Value = "Equipamiento y suministros médicos"
$ws.Range("E55").Value = 18866703.59299771
$ws.Range("F55").Value = 15664639472.27977
$ws.Range("G55").Value = ""
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 4

$ws.Range("A56").Value = 54
$ws.Range("B56").Value = 2023
$ws.Range("C56").Value = "Los Ríos"
$ws.Range("D56").Value = "Salud, servicios sanitarios y alimentación"
$ws.Range("E56").Value = 14083575.58454207
$ws.Range("F56").Value = 11703945947.56422
$ws.Range("G56").Value = ""
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 5

$ws.Range("A57").Value = 55
$ws.Range("B57").Value = 2023
$ws.Range("C57").Value = "Magallanes y Antártica"
$ws.Range("D57").Value = ""
$ws.Range("E57").Value = 61379845.62351756
$ws.Range("F57").Value = 53976679919
$ws.Range("G57").Value = 1488491.763933692
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 1

$ws.Range("A58").Value = 56
$ws.Range("B58").Value = 2023
$ws.Range("C58").Value = "Magallanes y Antártica"
$ws.Range("D58").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E58").Value = 33304843.31895265
$ws.Range("F58").Value = 27507200527.17492
$ws.Range("G58").Value = ""
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 2

$ws.Range("A59").Value = 57
$ws.Range("B59").Value = 2023
$ws.Range("C59").Value = "Magallanes y Antártica"
$ws.Range("D59").Value = "Obras MINVU"
$ws.Range("E59").Value = 22792563.78664332
$ws.Range("F59").Value = 19564617769.97734
$ws.Range("G59").Value = ""
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 3

$ws.Range("A60").Value = 58
$ws.Range("B60").Value = 2023
$ws.Range("C60").Value = "Magallanes y Antártica"
$ws.Range("D60").Value = "Equipamiento y suministros médicos"
$ws.Range("E60").Value = 16649682.38543603
$ws.Range("F60").Value = 14163299862.91663
$ws.Range("G60").Value = ""
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 4

$ws.Range("A61").Value = 59
$ws.Range("B61").Value = 2023
$ws.Range("C61").Value = "Magallanes y Antártica"
$ws.Range("D61").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E61").Value = 14623663.32306321
$ws.Range("F61").Value = 12259828719.89903
$ws.Range("G61").Value = ""
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 5

$ws.Range("A62").Value = 60
$ws.Range("B62").Value = 2023
$ws.Range("C62").Value = "Maule"
$ws.Range("D62").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E62").Value = 74252566.25951196
$ws.Range("F62").Value = 62616581647.55495
$ws.Range("G62").Value = ""
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 1

$ws.Range("A63").Value = 61
$ws.Range("B63").Value = 2023
$ws.Range("C63").Value = "Maule"
$ws.Range("D63").Value = ""
$ws.Range("E63").Value = 60146679.88510414
$ws.Range("F63").Value = 48828564123
$ws.Range("G63").Value = 1354604.918297481
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 2

$ws.Range("A64").Value = 62
$ws.Range("B64").Value = 2023
$ws.Range("C64").Value = "Maule"
$ws.Range("D64").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E64").Value = 56465735.1633732
$ws.Range("F64").Value = 47324148844.73659
$ws.Range("G64").Value = ""
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 3

$ws.Range("A65").Value = 63
$ws.Range("B65").Value = 2023
$ws.Range("C65").Value = "Maule"
$ws.Range("D65").Value = "Equipamiento y suministros médicos"
$ws.Range("E65").Value = 53470197.19490345
$ws.Range("F65").Value = 44827022219.4157
$ws.Range("G65").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 4

$ws.Range("A66").Value = 64
$ws.Range("B66").Value = 2023
$ws.Range("C66").Value = "Maule"
$ws.Range("D66").Value = "Obras MINVU"
$ws.Range("E66").Value = 43726748.89731067
$ws.Range("F66").Value = 37779561009.98734
$ws.Range("G66").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 5

$ws.Range("A67").Value = 65
$ws.Range("B67").Value = 2023
$ws.Range("C67").Value = "Metropolitana"
$ws.Range("D67").Value = "Medicamentos y productos farmacéuticos"
$ws.Range("E67").Value = 1245880871.338465
$ws.Range("F67").Value = 1031082771984.208
$ws.Range("G67").Value = ""
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 1

$ws.Range("A68").Value = 66
$ws.Range("B68").Value = 2023
$ws.Range("C68").Value = "Metropolitana"
$ws.Range("D68").Value = ""
$ws.Range("E68").Value = 512794630.8053578
$ws.Range("F68").Value = 417101963736
$ws.Range("G68").Value = 11669153.75548899
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 2

$ws.Range("A69").Value = 67
$ws.Range("B69").Value = 2023
$ws.Range("C69").Value = "Metropolitana"
$ws.Range("D69").Value = "Equipamiento y suministros médicos"
$ws.Range("E69").Value = 449948213.4182352
$ws.Range("F69").Value = 377021198025.541
$ws.Range("G69").Value = ""
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 3

$ws.Range("A70").Value = 68
$ws.Range("B70").Value = 2023
$ws.Range("C70").Value = "Metropolitana"
$ws.Range("D70").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E70").Value = 381535217.8998587
$ws.Range("F70").Value = 323310290829.6418
$ws.Range("G70").Value = ""
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 4

$ws.Range("A71").Value = 69
$ws.Range("B71").Value = 2023
$ws.Range("C71").Value = "Metropolitana"
$ws.Range("D71").Value = "Servicios basados en ingeniería, ciencias sociales y tecnología de la información"
$ws.Range("E71").Value = 351973947.2715064
$ws.Range("F71").Value = 296273919127.8158
$ws.Range("G71").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 5

$ws.Range("A72").Value = 70
$ws.Range("B72").Value = 2023
$ws.Range("C72").Value = "Ñuble"
$ws.Range("D72").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E72").Value = 123254360.5859251
$ws.Range("F72").Value = 101588784520.8411
$ws.Range("G72").Value = ""
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 1

$ws.Range("A73").Value = 71
$ws.Range("B73").Value = 2023
$ws.Range("C73").Value = "Ñuble"
$ws.Range("D73").Value = "Sin Información"
$ws.Range("E73").Value = 24443328.93593493
$ws.Range("F73").Value = 21140457873.18453
$ws.Range("G73").Value = 583912.405124515
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 2

$ws.Range("A74").Value = 72
$ws.Range("B74").Value = 2023
$ws.Range("C74").Value = "Ñuble"
$ws.Range("D74").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E74").Value = 22680111.59829282
$ws.Range("F74").Value = 19035314624.30814
$ws.Range("G74").Value = ""
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 3

$ws.Range("A75").Value = 73
$ws.Range("B75").Value = 2023
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E75").Value = 21972388.90992018
$ws.Range("F75").Value = 18453049421.75361
$ws.Range("G75").Value = ""
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 4

$ws.Range("A76").Value = 74
$ws.Range("B76").Value = 2023
$ws.Range("C76").Value = "Ñuble"
$ws.Range("D76").Value = "Equipamiento y suministros médicos"
$ws.Range("E76").Value = 19379773.17125701
$ws.Range("F76").Value = 16462359695.295
$ws.Range("G76").Value = ""
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 5

$ws.Range("A77").Value = 75
$ws.Range("B77").Value = 2023
$ws.Range("C77").Value = "Sin información"
$ws.Range("D77").Value = "Tecnologías de la información, telecomunicaciones y radiodifusión"
$ws.Range("E77").Value = 5053.726229795701
$ws.Range("F77").Value = 4411800
$ws.Range("G77").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 1

$ws.Range("A78").Value = 76
$ws.Range("B78").Value = 2023
$ws.Range("C78").Value = "Sin información"
$ws.Range("D78").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E78").Value = 4732.721525935397
$ws.Range("F78").Value = 3919818
$ws.Range("G78").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 2

$ws.Range("A79").Value = 77
$ws.Range("B79").Value = 2023
$ws.Range("C79").Value = "Sin información"
$ws.Range("D79").Value = "Vehículos y equipamiento en general"
$ws.Range("E79").Value = 4410.956829840198
$ws.Range("F79").Value = 4089041
$ws.Range("G79").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 3

$ws.Range("A80").Value = 78
$ws.Range("B80").Value = 2023
$ws.Range("C80").Value = "Sin información"
$ws.Range("D80").Value = "Medicamentos y productos farmacéuticos"
$ws.Range("E80").Value = 2655.174627710657
$ws.Range("F80").Value = 2312466
$ws.Range("G80").Value = ""
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 4

$ws.Range("A81").Value = 79
$ws.Range("B81").Value = 2023
$ws.Range("C81").Value = "Sin información"
$ws.Range("D81").Value = "Equipos, accesorios y suministros de oficina"
$ws.Range("E81").Value = 2426.65818972913
$ws.Range("F81").Value = 2108812
$ws.Range("G81").Value = ""
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 5

$ws.Range("A82").Value = 80
$ws.Range("B82").Value = 2023
$ws.Range("C82").Value = "Tarapacá"
$ws.Range("D82").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E82").Value = 29905064.68467002
$ws.Range("F82").Value = 25341561096.9814
$ws.Range("G82").Value = ""
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 1

$ws.Range("A83").Value = 81
$ws.Range("B83").Value = 2023
$ws.Range("C83").Value = "Tarapacá"
$ws.Range("D83").Value = "Equipamiento y suministros médicos"
$ws.Range("E83").Value = 20863069.25639435
$ws.Range("F83").Value = 17460996541.85113
$ws.Range("G83").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 2

$ws.Range("A84").Value = 82
$ws.Range("B84").Value = 2023
$ws.Range("C84").Value = "Tarapacá"
$ws.Range("D84").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E84").Value = 17363264.73087526
$ws.Range("F84").Value = 14783044416.01747
$ws.Range("G84").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 3

$ws.Range("A85").Value = 83
$ws.Range("B85").Value = 2023
$ws.Range("C85").Value = "Tarapacá"
$ws.Range("D85").Value = "Servicios de limpieza industrial"
$ws.Range("E85").Value = 13231942.77414693
$ws.Range("F85").Value = 11076150912.72366
$ws.Range("G85").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 4

$ws.Range("A86").Value = 84
$ws.Range("B86").Value = 2023
$ws.Range("C86").Value = "Tarapacá"
$ws.Range("D86").Value = "Servicios de transporte, almacenaje y correo"
$ws.Range("E86").Value = 10425432.63437072
$ws.Range("F86").Value = 8780605981.992975
$ws.Range("G86").Value = ""
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 5

$ws.Range("A87").Value = 85
$ws.Range("B87").Value = 2023
$ws.Range("C87").Value = "Valparaíso"
$ws.Range("D87").Value = "Servicios de construcción y mantenimiento"
$ws.Range("E87").Value = 114031351.5085096
$ws.Range("F87").Value = 95629652457.0411
$ws.Range("G87").Value = ""
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 1

$ws.Range("A88").Value = 86
$ws.Range("B88").Value = 2023
$ws.Range("C88").Value = "Valparaíso"
$ws.Range("D88").Value = "Equipamiento y suministros médicos"
$ws.Range("E88").Value = 106835103.1967579
$ws.Range("F88").Value = 89867292614.95744
$ws.Range("G88").Value = ""
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 2

$ws.Range("A89").Value = 87
$ws.Range("B89").Value = 2023
$ws.Range("C89").Value = "Valparaíso"
$ws.Range("D89").Value = ""
$ws.Range("E89").Value = 82064466.81275922
$ws.Range("F89").Value = 66280364260
$ws.Range("G89").Value = 1852515.13107252
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 3

$ws.Range("A90").Value = 88
$ws.Range("B90").Value = 2023
$ws.Range("C90").Value = "Valparaíso"
$ws.Range("D90").Value = "Servicios profesionales, administrativos y consultorías de gestión empresarial"
$ws.Range("E90").Value = 81165262.02298148
$ws.Range("F90").Value = 67372669016.7013
$ws.Range("G90").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 4

$ws.Range("A91").Value = 89
$ws.Range("B91").Value = 2023
$ws.Range("C91").Value = "Valparaíso"
$ws.Range("D91").Value = "Servicios de limpieza industrial"
$ws.Range("E91").Value = 70427948.90482001
$ws.Range("F91").Value = 58809972770.20901
$ws.Range("G91").Value = ""
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 5
# Apply column-A formatting (bold, centered, thin border) to the newly added rows 83-91,
# matching the style already used by column A in rows 2-82.
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83:A91").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
